$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegSpec_Example")

# capture original comment texts before row shift
$t_b12 = $ws.Range("B12").Comment.Text()
$t_d14 = $ws.Range("D14").Comment.Text()
$t_d17 = $ws.Range("D17").Comment.Text()

# delete B10 comment entirely
$ws.Range("B10").Comment.Delete()

# delete old comments before moving rows (they don't auto move)
$ws.Range("B12").Comment.Delete()
$ws.Range("D14").Comment.Delete()
$ws.Range("D17").Comment.Delete()

# Step 1: insert 2 new rows before row 7
$ws.Rows("7:8").Insert()

# Step 2: delete the old RegStructure row, now at row 12
$ws.Rows("12:12").Delete()

# re-add comments at shifted positions
$ws.Range("B13").AddComment($t_b12)
$ws.Range("D15").AddComment($t_d14)
$ws.Range("D18").AddComment($t_d17)

# Update RegFile table cell values
$ws.Range("B6").Value = "APB4"
$ws.Range("A7").Value = "DataWidth"
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = "Only 32, not changed"
$ws.Range("A8").Value = "AddrWidth"
$ws.Range("B8").Value = 16
$ws.Range("B10").Value = "No/Yes"
$ws.Range("B11").Value = "WProt/Sec"

# Update Register table header row values
$ws.Range("B13").Value = "Gen/-"
$ws.Range("E13").Value = "16'h0000"
